$d = $word.ActiveDocument

# =====================================================================
# PART 1: Insert a new "Meta description" paragraph right after the
# title paragraph (Heading1 "Play Blood Night Slot for Free - A
# Detailed Review"). The new paragraph contains an empty leading run,
# a bold "Meta description" run and a plain run with the rest of the
# text, and uses the default (Normal) paragraph style.
# =====================================================================

# --- Build the bold+plain formatted text fragment on a scratch
#     paragraph inserted after a plain, unformatted body paragraph so
#     that it doesn't pick up bold/italic/heading formatting. ---
$scratchAnchor = $d.Paragraphs(3)
$scratchInsertPoint = $scratchAnchor.Range.Duplicate
$scratchInsertPoint.Collapse(0)
$scratchInsertPoint.InsertParagraphAfter()
$scratchIdx = 4
$scratchPara = $d.Paragraphs($scratchIdx)
$scratchStart = $scratchPara.Range.Start
$fullText = "Meta description: Discover the features of Tuko's Blood Night slot game, including jackpots, graphics, symbols and Return to Player rate. Play for free and enjoy!"
$scratchPara.Range.InsertAfter($fullText)
$boldLen = "Meta description".Length
$boldRange = $d.Range($scratchStart, $scratchStart + $boldLen)
$boldRange.Font.Bold = 1

# copy the two-run range (bold run + plain run), excluding the
# paragraph mark
$copySrc = $d.Range($scratchStart, $scratchStart + $fullText.Length)
$copySrc.Copy()

# remove the scratch paragraph again (including the paragraph mark
# that was added before it)
$scratchPara2 = $d.Paragraphs($scratchIdx)
$delStart = $scratchPara2.Range.Start - 1
$delEnd = $scratchPara2.Range.End
$d.Range($delStart, $delEnd).Delete()

# --- Insert a new blank paragraph right after the title, switch it to
#     the default "Normal" style (so no pPr/pStyle is written), and
#     paste the previously copied bold+plain fragment into it. Pasting
#     into an untouched blank paragraph preserves its own leading
#     empty run. ---
$headingPara = $d.Paragraphs(2)
$headingPara.Range.InsertParagraphBefore()
$blankPara = $d.Paragraphs(2)
$blankPara.Style = "Normal"

$pasteRange = $blankPara.Range.Duplicate
$pasteRange.Collapse(0)
$pasteRange.Paste()

# =====================================================================
# PART 2: Near the end of the document, remove the duplicated bold
# "Play Blood Night Slot for Free - A Detailed Review" paragraph, and
# change the text of the final (italic) paragraph to the new image
# prompt text, keeping its italic formatting and leading empty run.
# =====================================================================

$count = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs($count - 1)
$delRange2 = $d.Range($boldTitlePara.Range.Start, $boldTitlePara.Range.End)
$delRange2.Delete()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$newImagePrompt = "Create a feature image for Blood Night that captures the game's dark and immersive atmosphere with a cartoon style. The image should prominently feature a Maya warrior wearing glasses with a happy expression on their face, set against a background of a cemetery at night with bats flying around. Try to incorporate some of the game's symbols, such as the gravestone, garlic cloves, and magic potions, into the image. Use a color scheme that matches the game's color palette, and make sure the imagery is visually striking and attention-grabbing."

$textRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$textRange.Text = $newImagePrompt
